$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from serial date 45791 to 45792 for rows 2-43
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45791) {
        $cell.Value = 45792
    }
}
